{"js": "const pairs = [\n  [\"2024-05-29 Wednesday\", \"2024-05-30 Thursday\"],\n  [\"9+33=42\", \"98-74=24\"],\n  [\"83-20=63\", \"68-17=51\"],\n  [\"51+32=83\", \"35+47=82\"],\n  [\"75-40=35\", \"51+28=79\"],\n  [\"3+1=4\", \"1+42=43\"],\n  [\"48-39=9\", \"66-54=12\"],\n  [\"50-14=36\", \"55-43=12\"],\n  [\"64+32=96\", \"41-38=3\"],\n  [\"16+47=63\", \"51-20=31\"],\n  [\"0+86=86\", \"33-7=26\"],\n  [\"32+60=92\", \"73-32=41\"],\n  [\"2+77=79\", \"15+40=55\"],\n  [\"92-40=52\", \"23+26=49\"],\n  [\"51-48=3\", \"5+43=48\"],\n  [\"79-75=4\", \"83-18=65\"],\n  [\"65-58=7\", \"12+71=83\"],\n  [\"75-15=60\", \"10+71=81\"],\n  [\"42+0=42\", \"16+21=37\"],\n  [\"71+9=80\", \"52-33=19\"],\n  [\"75+8=83\", \"53+37=90\"],\n  [\"40+29=69\", \"37+39=76\"],\n  [\"58-10=48\", \"50-43=7\"],\n  [\"28+22=50\", \"55+16=71\"],\n  [\"55-10=45\", \"57+11=68\"],\n  [\"47+7=54\", \"41+43=84\"],\n  [\"27+20=47\", \"86-3=83\"],\n  [\"14+12=26\", \"96-85=11\"],\n  [\"5+91=96\", \"9+54=63\"],\n  [\"68-12=56\", \"78-73=5\"],\n  [\"23+57=80\", \"45-33=12\"],\n  [\"34+27=61\", \"30+37=67\"],\n  [\"27+43=70\", \"22+33=55\"],\n  [\"53+9=62\", \"62-17=45\"],\n  [\"47-35=12\", \"76-0=76\"],\n  [\"61-48=13\", \"55-30=25\"],\n  [\"56-16=40\", \"43+52=95\"],\n  [\"68+27=95\", \"67-20=47\"],\n  [\"15+42=57\", \"57+26=83\"],\n  [\"42+25=67\", \"72-43=29\"],\n  [\"96-47=49\", \"45+40=85\"],\n  [\"65+12=77\", \"30+60=90\"],\n  [\"57-17=40\", \"44-30=14\"],\n  [\"44-9=35\", \"25-0=25\"],\n  [\"79-50=29\", \"47+5=52\"],\n  [\"33+39=72\", \"75-35=40\"],\n  [\"16+15=31\", \"64+10=74\"],\n  [\"99-84=15\", \"26-23=3\"],\n  [\"26-13=13\", \"32-25=7\"],\n  [\"57+42=99\", \"68+25=93\"],\n  [\"60-26=34\", \"21+63=84\"],\n  [\"33+4=37\", \"85-6=79\"],\n  [\"46+44=90\", \"42-31=11\"],\n  [\"6+71=77\", \"4+94=98\"],\n  [\"53-20=33\", \"36-6=30\"],\n  [\"72-59=13\", \"57+24=81\"],\n  [\"3+51=54\", \"5+89=94\"],\n  [\"52+29=81\", \"93-70=23\"],\n  [\"22+69=91\", \"92-70=22\"],\n  [\"18+61=79\", \"33+13=46\"],\n  [\"61+17=78\", \"7-6=1\"],\n  [\"7+27=34\", \"80-63=17\"],\n  [\"73+0=73\", \"86-65=21\"],\n  [\"87-40=47\", \"71-54=17\"],\n  [\"56+29=85\", \"69-37=32\"],\n  [\"12+43=55\", \"26-19=7\"],\n  [\"32-11=21\", \"50-17=33\"],\n  [\"21+12=33\", \"18+52=70\"],\n  [\"15+1=16\", \"48-9=39\"],\n  [\"56+40=96\", \"35+44=79\"],\n  [\"72-28=44\", \"19+42=61\"],\n  [\"13+4=17\", \"2+20=22\"],\n  [\"28+65=93\", \"97-80=17\"],\n  [\"36+63=99\", \"74-10=64\"],\n  [\"91-53=38\", \"39-14=25\"],\n  [\"78-40=38\", \"28+70=98\"],\n  [\"54-27=27\", \"62-49=13\"],\n  [\"0+58=58\", \"86-52=34\"],\n  [\"41-2=39\", \"46-22=24\"],\n  [\"70-69=1\", \"25+64=89\"],\n  [\"93-76=17\", \"73+9=82\"],\n  [\"87-71=16\", \"43+30=73\"],\n  [\"60+12=72\", \"33+8=41\"],\n  [\"20+13=33\", \"43-1=42\"],\n  [\"44+9=53\", \"8+13=21\"],\n  [\"72-26=46\", \"58-27=31\"],\n  [\"18+20=38\", \"31-18=13\"],\n  [\"19+0=19\", \"51-12=39\"],\n  [\"79+17=96\", \"31+60=91\"],\n  [\"60+23=83\", \"16-8=8\"],\n  [\"19+32=51\", \"24-2=22\"],\n  [\"18+45=63\", \"24+29=53\"],\n  [\"53-30=23\", \"49-40=9\"],\n  [\"54-48=6\", \"98-8=90\"],\n  [\"80-41=39\", \"77-58=19\"],\n  [\"93-56=37\", \"17+67=84\"],\n  [\"96+1=97\", \"16+58=74\"],\n  [\"10+70=80\", \"97-12=85\"],\n  [\"65-50=15\", \"52-4=48\"],\n  [\"85+11=96\", \"83-56=27\"],\n  [\"89-29=60\", \"47-19=28\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"2024-05-29 Wednesday\", \"2024-05-30 Thursday\"),\n  @(\"9+33=42\", \"98-74=24\"),\n  @(\"83-20=63\", \"68-17=51\"),\n  @(\"51+32=83\", \"35+47=82\"),\n  @(\"75-40=35\", \"51+28=79\"),\n  @(\"3+1=4\", \"1+42=43\"),\n  @(\"48-39=9\", \"66-54=12\"),\n  @(\"50-14=36\", \"55-43=12\"),\n  @(\"64+32=96\", \"41-38=3\"),\n  @(\"16+47=63\", \"51-20=31\"),\n  @(\"0+86=86\", \"33-7=26\"),\n  @(\"32+60=92\", \"73-32=41\"),\n  @(\"2+77=79\", \"15+40=55\"),\n  @(\"92-40=52\", \"23+26=49\"),\n  @(\"51-48=3\", \"5+43=48\"),\n  @(\"79-75=4\", \"83-18=65\"),\n  @(\"65-58=7\", \"12+71=83\"),\n  @(\"75-15=60\", \"10+71=81\"),\n  @(\"42+0=42\", \"16+21=37\"),\n  @(\"71+9=80\", \"52-33=19\"),\n  @(\"75+8=83\", \"53+37=90\"),\n  @(\"40+29=69\", \"37+39=76\"),\n  @(\"58-10=48\", \"50-43=7\"),\n  @(\"28+22=50\", \"55+16=71\"),\n  @(\"55-10=45\", \"57+11=68\"),\n  @(\"47+7=54\", \"41+43=84\"),\n  @(\"27+20=47\", \"86-3=83\"),\n  @(\"14+12=26\", \"96-85=11\"),\n  @(\"5+91=96\", \"9+54=63\"),\n  @(\"68-12=56\", \"78-73=5\"),\n  @(\"23+57=80\", \"45-33=12\"),\n  @(\"34+27=61\", \"30+37=67\"),\n  @(\"27+43=70\", \"22+33=55\"),\n  @(\"53+9=62\", \"62-17=45\"),\n  @(\"47-35=12\", \"76-0=76\"),\n  @(\"61-48=13\", \"55-30=25\"),\n  @(\"56-16=40\", \"43+52=95\"),\n  @(\"68+27=95\", \"67-20=47\"),\n  @(\"15+42=57\", \"57+26=83\"),\n  @(\"42+25=67\", \"72-43=29\"),\n  @(\"96-47=49\", \"45+40=85\"),\n  @(\"65+12=77\", \"30+60=90\"),\n  @(\"57-17=40\", \"44-30=14\"),\n  @(\"44-9=35\", \"25-0=25\"),\n  @(\"79-50=29\", \"47+5=52\"),\n  @(\"33+39=72\", \"75-35=40\"),\n  @(\"16+15=31\", \"64+10=74\"),\n  @(\"99-84=15\", \"26-23=3\"),\n  @(\"26-13=13\", \"32-25=7\"),\n  @(\"57+42=99\", \"68+25=93\"),\n  @(\"60-26=34\", \"21+63=84\"),\n  @(\"33+4=37\", \"85-6=79\"),\n  @(\"46+44=90\", \"42-31=11\"),\n  @(\"6+71=77\", \"4+94=98\"),\n  @(\"53-20=33\", \"36-6=30\"),\n  @(\"72-59=13\", \"57+24=81\"),\n  @(\"3+51=54\", \"5+89=94\"),\n  @(\"52+29=81\", \"93-70=23\"),\n  @(\"22+69=91\", \"92-70=22\"),\n  @(\"18+61=79\", \"33+13=46\"),\n  @(\"61+17=78\", \"7-6=1\"),\n  @(\"7+27=34\", \"80-63=17\"),\n  @(\"73+0=73\", \"86-65=21\"),\n  @(\"87-40=47\", \"71-54=17\"),\n  @(\"56+29=85\", \"69-37=32\"),\n  @(\"12+43=55\", \"26-19=7\"),\n  @(\"32-11=21\", \"50-17=33\"),\n  @(\"21+12=33\", \"18+52=70\"),\n  @(\"15+1=16\", \"48-9=39\"),\n  @(\"56+40=96\", \"35+44=79\"),\n  @(\"72-28=44\", \"19+42=61\"),\n  @(\"13+4=17\", \"2+20=22\"),\n  @(\"28+65=93\", \"97-80=17\"),\n  @(\"36+63=99\", \"74-10=64\"),\n  @(\"91-53=38\", \"39-14=25\"),\n  @(\"78-40=38\", \"28+70=98\"),\n  @(\"54-27=27\", \"62-49=13\"),\n  @(\"0+58=58\", \"86-52=34\"),\n  @(\"41-2=39\", \"46-22=24\"),\n  @(\"70-69=1\", \"25+64=89\"),\n  @(\"93-76=17\", \"73+9=82\"),\n  @(\"87-71=16\", \"43+30=73\"),\n  @(\"60+12=72\", \"33+8=41\"),\n  @(\"20+13=33\", \"43-1=42\"),\n  @(\"44+9=53\", \"8+13=21\"),\n  @(\"72-26=46\", \"58-27=31\"),\n  @(\"18+20=38\", \"31-18=13\"),\n  @(\"19+0=19\", \"51-12=39\"),\n  @(\"79+17=96\", \"31+60=91\"),\n  @(\"60+23=83\", \"16-8=8\"),\n  @(\"19+32=51\", \"24-2=22\"),\n  @(\"18+45=63\", \"24+29=53\"),\n  @(\"53-30=23\", \"49-40=9\"),\n  @(\"54-48=6\", \"98-8=90\"),\n  @(\"80-41=39\", \"77-58=19\"),\n  @(\"93-56=37\", \"17+67=84\"),\n  @(\"96+1=97\", \"16+58=74\"),\n  @(\"10+70=80\", \"97-12=85\"),\n  @(\"65-50=15\", \"52-4=48\"),\n  @(\"85+11=96\", \"83-56=27\"),\n  @(\"89-29=60\", \"47-19=28\")\n)\n\nforeach ($pair in $pairs) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        throw \"No match found for: $oldText\"\n    }\n}\n"}
